$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Inha"
$ws.Range("C2").Value = "Tgfbr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3688166666666666
$ws.Range("H2").Value = 1.10645
$ws.Range("I2").Value = 0.480574717760894
$ws.Range("J2").Value = 0.5812056951802134
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 35.789624
$ws.Range("N2").Value = 71.57924800000001
$ws.Range("O2").Value = 0.258139457682779
$ws.Range("P2").Value = 0.1993778771086309
$ws.Range("Q2").Value = 13.19980982493333
$ws.Range("R2").Value = 79.19885894960001
$ws.Range("S2").Value = 0.1240552970188517
$ws.Range("T2").Value = 0.115879557668477

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Inha"
$ws.Range("C3").Value = "Tgfbr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3688166666666666
$ws.Range("H3").Value = 1.10645
$ws.Range("I3").Value = 0.480574717760894
$ws.Range("J3").Value = 0.5812056951802134
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 60.113367
$ws.Range("N3").Value = 180.340101
$ws.Range("O3").Value = 0.4335790718803266
$ws.Range("P3").Value = 0.5023219368682956
$ws.Range("Q3").Value = 22.17081163905
$ws.Range("R3").Value = 199.53730475145
$ws.Range("S3").Value = 0.2083671400959183
$ws.Range("T3").Value = 0.291952370521809

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Inha"
$ws.Range("C4").Value = "Tgfbr3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3688166666666666
$ws.Range("H4").Value = 1.10645
$ws.Range("I4").Value = 0.480574717760894
$ws.Range("J4").Value = 0.5812056951802134
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2072186666666667
$ws.Range("N4").Value = 0.621656
$ws.Range("O4").Value = 0.001494603973349423
$ws.Range("P4").Value = 0.001731569652308208
$ws.Range("Q4").Value = 0.0764256979111111
$ws.Range("R4").Value = 0.6878312811999999
$ws.Range("S4").Value = 0.0007182688826567098
$ws.Range("T4").Value = 0.001006398143522753

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Inha"
$ws.Range("C5").Value = "Tgfbr3"
$ws.Range("D5").Value = "Neutro"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3688166666666666
$ws.Range("H5").Value = 1.10645
$ws.Range("I5").Value = 0.480574717760894
$ws.Range("J5").Value = 0.5812056951802134
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 21.40334366666667
$ws.Range("N5").Value = 64.210031
$ws.Range("O5").Value = 0.1543756795743782
$ws.Range("P5").Value = 0.178851553034748
$ws.Range("Q5").Value = 7.893909866661111
$ws.Range("R5").Value = 71.04518879995
$ws.Range("S5").Value = 0.07418904864060301
$ws.Range("T5").Value = 0.1039495412156215

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Inha"
$ws.Range("C6").Value = "Tgfbr3"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.3688166666666666
$ws.Range("H6").Value = 1.10645
$ws.Range("I6").Value = 0.480574717760894
$ws.Range("J6").Value = 0.5812056951802134
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 21.1309775
$ws.Range("N6").Value = 42.261955
$ws.Range("O6").Value = 0.1524111868891667
$ws.Range("P6").Value = 0.1177170633360173
$ws.Range("Q6").Value = 7.793456684958333
$ws.Range("R6").Value = 46.76074010975
$ws.Range("S6").Value = 0.07324496312286413
$ws.Range("T6").Value = 0.06841782763078313

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Inha"
$ws.Range("C7").Value = "Tgfbr3"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.3986325
$ws.Range("H7").Value = 0.797265
$ws.Range("I7").Value = 0.519425282239106
$ws.Range("J7").Value = 0.4187943048197866
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 35.789624
$ws.Range("N7").Value = 71.57924800000001
$ws.Range("O7").Value = 0.258139457682779
$ws.Range("P7").Value = 0.1993778771086309
$ws.Range("Q7").Value = 14.26690728918
$ws.Range("R7").Value = 57.06762915672
$ws.Range("S7").Value = 0.1340841606639272
$ws.Range("T7").Value = 0.08349831944015391

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Inha"
$ws.Range("C8").Value = "Tgfbr3"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.3986325
$ws.Range("H8").Value = 0.797265
$ws.Range("I8").Value = 0.519425282239106
$ws.Range("J8").Value = 0.4187943048197866
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 60.113367
$ws.Range("N8").Value = 180.340101
$ws.Range("O8").Value = 0.4335790718803266
$ws.Range("P8").Value = 0.5023219368682956
$ws.Range("Q8").Value = 23.9631417706275
$ws.Range("R8").Value = 143.778850623765
$ws.Range("S8").Value = 0.2252119317844083
$ws.Range("T8").Value = 0.2103695663464866

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Inha"
$ws.Range("C9").Value = "Tgfbr3"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.3986325
$ws.Range("H9").Value = 0.797265
$ws.Range("I9").Value = 0.519425282239106
$ws.Range("J9").Value = 0.4187943048197866
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.2072186666666667
$ws.Range("N9").Value = 0.621656
$ws.Range("O9").Value = 0.001494603973349423
$ws.Range("P9").Value = 0.001731569652308208
$ws.Range("Q9").Value = 0.08260409514
$ws.Range("R9").Value = 0.49562457084
$ws.Range("S9").Value = 0.0007763350906927134
$ws.Range("T9").Value = 0.0007251715087854556

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Inha"
$ws.Range("C10").Value = "Tgfbr3"
$ws.Range("D10").Value = "Neutro"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.3986325
$ws.Range("H10").Value = 0.797265
$ws.Range("I10").Value = 0.519425282239106
$ws.Range("J10").Value = 0.4187943048197866
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 21.40334366666667
$ws.Range("N10").Value = 64.210031
$ws.Range("O10").Value = 0.1543756795743782
$ws.Range("P10").Value = 0.178851553034748
$ws.Range("Q10").Value = 8.532068394202501
$ws.Range("R10").Value = 51.192410365215
$ws.Range("S10").Value = 0.08018663093377519
$ws.Range("T10").Value = 0.07490201181912647

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Inha"
$ws.Range("C11").Value = "Tgfbr3"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.3986325
$ws.Range("H11").Value = 0.797265
$ws.Range("I11").Value = 0.519425282239106
$ws.Range("J11").Value = 0.4187943048197866
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 21.1309775
$ws.Range("N11").Value = 42.261955
$ws.Range("O11").Value = 0.1524111868891667
$ws.Range("P11").Value = 0.1177170633360173
$ws.Range("Q11").Value = 8.42349438826875
$ws.Range("R11").Value = 33.693977553075
$ws.Range("S11").Value = 0.07916622376630253
$ws.Range("T11").Value = 0.04929923570523415
